$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new metric data row (row 13), matching the existing pattern of
# timestamp strings in column A and numeric metric values in column B.
$ws.Cells.Item(13, 1).Value = "2025-04-28 11:57:10"
$ws.Cells.Item(13, 2).Value = 239
